$wb = $excel.ActiveWorkbook

# 1. Instructions sheet: trim trailing space from B15
$wsInstructions = $wb.Worksheets.Item("Instructions")
$wsInstructions.Unprotect()
$wsInstructions.Range("B15").Value = "- Affinity: Spike protein binding affinity; inhibition of ACE2 binding; ELISA for Spike"
$wsInstructions.Protect()

# 2. Antibodies sheet: narrow the Light chain data validation list from $C$2:$C$4 to $C$2:$C$3
$wsAntibodies = $wb.Worksheets.Item("Antibodies")
$dRange = $wsAntibodies.Range("D2:D100")
$dRange.Validation.Modify(3, 1, 1, "=Terminology!`$C`$2:`$C`$3")

# 3. Terminology sheet: clear the "unknown" light-chain value in C4
$wsTerminology = $wb.Worksheets.Item("Terminology")
$wsTerminology.Unprotect()
$wsTerminology.Range("C4").ClearContents()
$wsTerminology.Protect()
